$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (B..L) right
$ws.Range("A1").EntireColumn.Insert()

# Populate the new "id" column
$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = "eedfd009-a800-4426-95b8-1bc2251dbeb9"

# Match the formatting of the adjacent (now-shifted) header/data cells
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set column A width to fit the new UUID values (closest representable width
# to the 40.7109375 "bestFit" value Excel's real font metrics produced)
$ws.Range("A1").ColumnWidth = 39.75

# Update selection to mirror the saved workbook view state
$ws.Range("A3").Select()
